$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.799.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.543.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.84%  "
$ws.Range("E7").Value = "  -0.65%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.934.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.618.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.817"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.811.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("E20").Value = "  -0.53%  "
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.87%  "
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.42%  "
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.84%  "
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("E36").Value = "  -3.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.25%  "
$ws.Range("E38").Value = "  -5.20%  "
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.117"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.33%  "
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0299"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.997.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.788.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "79.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.33%  "
